# Item Data workbook maintenance edit:
#  - "Burger" (row 4) is now flagged as a Special Item (D4: FALSE -> TRUE)
#  - "test" item #6 (row 7) is no longer Active (E7: TRUE -> FALSE)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = Burger -> mark "Special Item" (column D) as TRUE
$ws.Range("D4").Value = $true

# Row 7 = test (Item ID 6) -> mark "Is Active" (column E) as FALSE
$ws.Range("E7").Value = $false
